# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape (numbers shifted by +/-1).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2282
    $ws.Range("F5").Value = 1094
    $ws.Range("F6").Value = 864
    $ws.Range("F8").Value = 5852
}
